$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Author name
$ws.Range("C3").Value = "Ridham Sood"

# Column G (Expected Result) for accessor-style rows first
$ws.Range("G7").Value = "Set attributes to the input values."
$ws.Range("G12").Value = "Accessor return client_number attribute."
$ws.Range("G13").Value = "Accessor return first_name attribute."
$ws.Range("G14").Value = "Accessor return last_name attribute."
$ws.Range("G15").Value = "Accessor return email_address attribute."

# Column F (Method Inputs) - common constructor call used across several rows
$ws.Range("F7").Value = 'client = Client(1, "Joe", "Henderson", "joehenderson9@gmail.com")'
$ws.Range("F12").Value = 'client = Client(1, "Joe", "Henderson", "joehenderson9@gmail.com")'
$ws.Range("F13").Value = 'client = Client(1, "Joe", "Henderson", "joehenderson9@gmail.com")'
$ws.Range("F14").Value = 'client = Client(1, "Joe", "Henderson", "joehenderson9@gmail.com")'
$ws.Range("F15").Value = 'client = Client(1, "Joe", "Henderson", "joehenderson9@gmail.com")'
$ws.Range("F16").Value = 'client = Client(1, "Joe", "Henderson", "joehenderson9@gmail.com")'

$ws.Range("G16").Value = "Str returns valid statement."

$ws.Range("E12").Value = "from client.client import Client.           Import unittest"

# Column F for the invalid-input rows (8-11)
$ws.Range("F8").Value = 'client = Client("one", "Joe", "Henderson", "joehenderson9@gmail.com")'
$ws.Range("F9").Value = 'client = Client(1, "", "Henderson", "joehenderson9@gmail.com")'
$ws.Range("F10").Value = 'client = Client(1, "Joe", "", "joehenderson9@gmail.com")'
$ws.Range("F11").Value = 'client = Client(1, "Joe", "Henderson", "joehenderson9")'

# Column G for the invalid-input rows (11 -> 8, reverse order)
$ws.Range("G11").Value = "ValueError - Email address should be in the correct format."
$ws.Range("G10").Value = "ValueError - Last name cannot be blank."
$ws.Range("G9").Value = "ValueError - First name cannot be blank."
$ws.Range("G8").Value = "ValueError - Client number must be an int type."

$ws.Range("E16").Value = "from client.client import Client        Import unittest"

# Column E (Preconditions) - common import statement used across several rows
$ws.Range("E7").Value = "from client.client import Client           Import unittest"
$ws.Range("E8").Value = "from client.client import Client           Import unittest"
$ws.Range("E9").Value = "from client.client import Client           Import unittest"
$ws.Range("E10").Value = "from client.client import Client           Import unittest"
$ws.Range("E11").Value = "from client.client import Client           Import unittest"
$ws.Range("E13").Value = "from client.client import Client           Import unittest"
$ws.Range("E14").Value = "from client.client import Client           Import unittest"
$ws.Range("E15").Value = "from client.client import Client           Import unittest"

$ws.Range("E7").Select() | Out-Null

# Row heights recalculated by Excel's autofit/re-render after the content edits
$ws.Rows.Item(12).RowHeight = 57
$ws.Rows.Item(13).RowHeight = 59.4
$ws.Rows.Item(14).RowHeight = 58.2
$ws.Rows.Item(15).RowHeight = 61.2
$ws.Rows.Item(16).RowHeight = 62.4
$ws.Rows.Item(17).RowHeight = 31.2
$ws.Rows.Item(18).RowHeight = 31.2
$ws.Rows.Item(19).RowHeight = 31.2
$ws.Rows.Item(20).RowHeight = 31.2
$ws.Rows.Item(21).RowHeight = 31.2
$ws.Rows.Item(22).RowHeight = 31.2
$ws.Rows.Item(23).RowHeight = 31.2
$ws.Rows.Item(24).RowHeight = 31.2
$ws.Rows.Item(25).RowHeight = 31.2
$ws.Rows.Item(26).RowHeight = 31.2
$ws.Rows.Item(27).RowHeight = 31.2
$ws.Rows.Item(28).RowHeight = 31.2
$ws.Rows.Item(2).RowHeight = 73.2
$null
